$wb = $excel.ActiveWorkbook

# Delete row 16 ("Sheet" summary row) on the "optimization_parameters" sheet.
# This shifts the old row 17 up to row 16 and removes the now-unused
# shared string "Sheet" and its associated number format style.
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Activate()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(16).EntireRow.Select()

# Make "threshold_b" the active / selected sheet.
$ws2 = $wb.Worksheets.Item("threshold_b")
$ws2.Activate()
